$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424251357742037
$ws.Range("D2").Value = 0.3327257862217721
$ws.Range("E2").Value = 0.07331599993843696
$ws.Range("F2").Value = 8.8416736834983
$ws.Range("G2").Value = 0.002693157817373358
$ws.Range("J2").Value = 0.06846980932727931
$ws.Range("M2").Value = 3.704961644579555
$ws.Range("N2").Value = 1.422206768421688
$ws.Range("B3").Value = 0.132937712117851
$ws.Range("D3").Value = 0.2976225236081973
$ws.Range("E3").Value = 0.06373473003400676
$ws.Range("F3").Value = 8.612767752135028
$ws.Range("G3").Value = 0.002710927555352015
$ws.Range("J3").Value = 0.06758737773010282
$ws.Range("M3").Value = 3.48253898129667
$ws.Range("N3").Value = 1.431318878071025
$ws.Range("B4").Value = 0.1271826007833425
$ws.Range("D4").Value = 0.2764085397485019
$ws.Range("E4").Value = 0.05786952203919071
$ws.Range("F4").Value = 8.481048648594168
$ws.Range("G4").Value = 0.002722338300031042
$ws.Range("J4").Value = 0.06705173993185198
$ws.Range("M4").Value = 3.348696403160773
$ws.Range("N4").Value = 1.437671333302944
$ws.Range("B5").Value = 0.124855124924494
$ws.Range("D5").Value = 0.267843439248594
$ws.Range("E5").Value = 0.05548279524793998
$ws.Range("F5").Value = 8.429534558809394
$ws.Range("G5").Value = 0.002727115009620664
$ws.Range("J5").Value = 0.0668349815686895
$ws.Range("M5").Value = 3.294815736227349
$ws.Range("N5").Value = 1.440451369632171
$ws.Range("B6").Value = 0.124469726868071
$ws.Range("D6").Value = 0.2664258683314529
$ws.Range("E6").Value = 0.05508665560472537
$ws.Range("F6").Value = 8.421109642003415
$ws.Range("G6").Value = 0.002727915861971015
$ws.Range("J6").Value = 0.06679907972285548
$ws.Range("M6").Value = 3.285908144021306
$ws.Range("N6").Value = 1.440924575709843
$ws.Range("B7").Value = 0.127151139475842
$ws.Range("D7").Value = 0.2762927121004566
$ws.Range("E7").Value = 0.05783732144906395
$ws.Range("F7").Value = 8.48034523053812
$ws.Range("G7").Value = 0.00272240220610917
$ws.Range("J7").Value = 0.06704881055007661
$ws.Range("M7").Value = 3.34796710339711
$ws.Range("N7").Value = 1.437708049881479
$ws.Range("B8").Value = 0.1391393718328402
$ws.Range("D8").Value = 0.3205482549685996
$ws.Range("E8").Value = 0.07000804245964076
$ws.Range("F8").Value = 8.760878663305903
$ws.Range("G8").Value = 0.002699181643825504
$ws.Range("J8").Value = 0.06816424231885065
$ws.Range("M8").Value = 3.627689644949982
$ws.Range("N8").Value = 1.425191896030867
$ws.Range("B9").Value = 0.1632013034624151
$ws.Range("D9").Value = 0.4102941477815421
$ws.Range("E9").Value = 0.09406638417285507
$ws.Range("F9").Value = 9.383815710844544
$ws.Range("G9").Value = 0.002657567578034588
$ws.Range("J9").Value = 0.07040244611254565
$ws.Range("M9").Value = 4.199020113227107
$ws.Range("N9").Value = 1.406622458071354
$ws.Range("B10").Value = 0.1812135996092934
$ws.Range("D10").Value = 0.4784258943271311
$ws.Range("E10").Value = 0.1119338434516024
$ws.Range("F10").Value = 9.88990247004358
$ws.Range("G10").Value = 0.002629317983139971
$ws.Range("J10").Value = 0.07208072087731665
$ws.Range("M10").Value = 4.634409935311652
$ws.Range("N10").Value = 1.396574814320218
$ws.Range("B11").Value = 0.1894799153733544
$ws.Range("D11").Value = 0.5099878824534585
$ws.Range("E11").Value = 0.1201212641504483
$ws.Range("F11").Value = 10.13158132969437
$ws.Range("G11").Value = 0.002616956332876233
$ws.Range("J11").Value = 0.07285225769906845
$ws.Range("M11").Value = 4.836280531111555
$ws.Range("N11").Value = 1.39277496353364
$ws.Range("B12").Value = 0.1926204914766458
$ws.Range("D12").Value = 0.5220291009986795
$ws.Range("E12").Value = 0.1232316602000978
$ws.Range("F12").Value = 10.22482696967745
$ws.Range("G12").Value = 0.002612344417601091
$ws.Range("J12").Value = 0.07314563869430657
$ws.Range("M12").Value = 4.91330616976191
$ws.Range("N12").Value = 1.391446125285739
$ws.Range("B13").Value = 0.191943655814228
$ws.Range("D13").Value = 0.5194317083671081
$ws.Range("E13").Value = 0.1225613121668019
$ws.Range("F13").Value = 10.20466672876131
$ws.Range("G13").Value = 0.00261333461867544
$ws.Range("J13").Value = 0.07308239881773204
$ws.Range("M13").Value = 4.896690893691698
$ws.Range("N13").Value = 1.391727430508411
$ws.Range("B14").Value = 0.1897380866336817
$ws.Range("D14").Value = 0.5109766843898456
$ws.Range("E14").Value = 0.1203769493473885
$ws.Range("F14").Value = 10.13921762064672
$ws.Range("G14").Value = 0.00261657552850307
$ws.Range("J14").Value = 0.07287636961720523
$ws.Range("M14").Value = 4.842605619702937
$ws.Range("N14").Value = 1.392663437339806
$ws.Range("B15").Value = 0.1883884501110771
$ws.Range("D15").Value = 0.5058096108169252
$ws.Range("E15").Value = 0.1190403104327658
$ws.Range("F15").Value = 10.09935550253442
$ws.Range("G15").Value = 0.002618569651780974
$ws.Range("J15").Value = 0.07275033093217331
$ws.Range("M15").Value = 4.809553608242624
$ws.Range("N15").Value = 1.393251082015894
$ws.Range("B16").Value = 0.1806748241120602
$ws.Range("D16").Value = 0.4763753273219891
$ws.Range("E16").Value = 0.1114000944478661
$ws.Range("F16").Value = 9.874345447166206
$ws.Range("G16").Value = 0.002630135590863058
$ws.Range("J16").Value = 0.07203046632677967
$ws.Range("M16").Value = 4.621296723365475
$ws.Range("N16").Value = 1.396838589503588
$ws.Range("B17").Value = 0.1759612335246175
$ws.Range("D17").Value = 0.4584692638042043
$ws.Range("E17").Value = 0.106729299005508
$ws.Range("F17").Value = 9.739298905143698
$ws.Range("G17").Value = 0.002637355397655459
$ws.Range("J17").Value = 0.07159096012177102
$ws.Range("M17").Value = 4.506806418967926
$ws.Range("N17").Value = 1.399236297631759
$ws.Range("B18").Value = 0.1732569250608123
$ws.Range("D18").Value = 0.4482232019738319
$ws.Range("E18").Value = 0.1040482829085647
$ws.Range("F18").Value = 9.662697137341752
$ws.Range("G18").Value = 0.002641554161599876
$ws.Range("J18").Value = 0.07133892490685767
$ws.Range("M18").Value = 4.441311066693714
$ws.Range("N18").Value = 1.400687999156403
$ws.Range("B19").Value = 0.1723424682823946
$ws.Range("D19").Value = 0.4447629437124192
$ws.Range("E19").Value = 0.1031414388540952
$ws.Range("F19").Value = 9.636943100681833
$ws.Range("G19").Value = 0.002642983747182352
$ws.Range("J19").Value = 0.07125371855212137
$ws.Range("M19").Value = 4.419195685261144
$ws.Range("N19").Value = 1.401192013630606
$ws.Range("B20").Value = 0.1764622978525381
$ws.Range("D20").Value = 0.4603698521725619
$ws.Range("E20").Value = 0.1072259334476655
$ws.Range("F20").Value = 9.753563180918093
$ws.Range("G20").Value = 0.002636582072017113
$ws.Range("J20").Value = 0.0716376675639836
$ws.Range("M20").Value = 4.51895697789422
$ws.Range("N20").Value = 1.398973548954956
$ws.Range("B21").Value = 0.1903856370652051
$ws.Range("D21").Value = 0.5134576414493495
$ws.Range("E21").Value = 0.1210182663707897
$ws.Range("F21").Value = 10.15839408280408
$ws.Range("G21").Value = 0.002615621726921121
$ws.Range("J21").Value = 0.07293685190896326
$ws.Range("M21").Value = 4.858475711700436
$ws.Range("N21").Value = 1.392385527884628
$ws.Range("B22").Value = 0.1995453408346179
$ws.Range("D22").Value = 0.5486783719229607
$ws.Range("E22").Value = 0.130091356370464
$ws.Range("F22").Value = 10.43308126493912
$ws.Range("G22").Value = 0.002602325549161344
$ws.Range("J22").Value = 0.07379306451012724
$ws.Range("M22").Value = 5.083778058100648
$ws.Range("N22").Value = 1.388721174505662
$ws.Range("B23").Value = 0.1946511813049625
$ws.Range("D23").Value = 0.5298297636836082
$ws.Range("E23").Value = 0.1252429785648417
$ws.Range("F23").Value = 10.28552340637503
$ws.Range("G23").Value = 0.002609385529287689
$ws.Range("J23").Value = 0.07333541703989255
$ws.Range("M23").Value = 4.963206523029072
$ws.Range("N23").Value = 1.390618480336514
$ws.Range("B24").Value = 0.1762357492260236
$ws.Range("D24").Value = 0.459510446073466
$ws.Range("E24").Value = 0.1070013917739985
$ws.Range("F24").Value = 9.747111073771805
$ws.Range("G24").Value = 0.002636931542989832
$ws.Range("J24").Value = 0.07161654913130988
$ws.Range("M24").Value = 4.513462695299921
$ws.Range("N24").Value = 1.399092109487015
$ws.Range("B25").Value = 0.156633046539099
$ws.Range("D25").Value = 0.3856592775103991
$ws.Range("E25").Value = 0.08752991909873487
$ws.Range("F25").Value = 9.207095282009277
$ws.Range("G25").Value = 0.002668412221411994
$ws.Range("J25").Value = 0.06979122420657546
$ws.Range("M25").Value = 4.041845588682804
$ws.Range("N25").Value = 1.411011272598643
